# Append two new daily rows (2025-10-05 / serial 45935) to Sheet1:
#   row 70 -> 四方坪站充电量(kw)
#   row 71 -> 高岭站充电量(kw)
# mirroring the existing per-station daily rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

$row70 = @(
    45935,
    "四方坪站充电量(kw)",
    591.58000000000004,
    1127.3609999999996,
    586.33300000000008,
    369.3,
    508.61100000000016,
    637.76800000000003,
    354.12500000000006,
    207.398,
    139.733,
    150.70600000000002,
    260.17500000000001,
    233.191,
    871.71399999999994,
    935.64600000000007,
    548.70699999999999,
    550.36700000000008,
    463.19,
    172.41500000000002,
    113.03,
    89.75,
    96.320000000000007,
    62.82,
    50.09,
    13.46
)

$row71 = @(
    45935,
    "高岭站充电量(kw)",
    479.173,
    510.654,
    199.649,
    149.22800000000001,
    0,
    161.69499999999999,
    49.418000000000006,
    176.91300000000001,
    199.00700000000001,
    130.27699999999999,
    191.07499999999999,
    190.73,
    273.36400000000003,
    328.447,
    36.871000000000002,
    257.50799999999998,
    145.03800000000001,
    23.234000000000002,
    10.338999999999999,
    41.927999999999997,
    22,
    55.576999999999998,
    14.25,
    0
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "70").Value = $row70[$i]
    $ws.Range($cols[$i] + "71").Value = $row71[$i]
}

# Keep the workbook's selection state consistent with the extra two rows.
$ws.Range("F75").Select()
